$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 0.3484863333333333
$ws.Cells.Item(2, 8).Value = 1.045459
$ws.Cells.Item(2, 9).Value = 0.08515845388213966
$ws.Cells.Item(2, 10).Value = 0.08515845388213968
$ws.Cells.Item(2, 13).Value = 3.339352
$ws.Cells.Item(2, 14).Value = 10.018056
$ws.Cells.Item(2, 15).Value = 0.6054960700393903
$ws.Cells.Item(2, 16).Value = 0.6054960700393903
$ws.Cells.Item(2, 17).Value = 1.163718534189333
$ws.Cells.Item(2, 18).Value = 10.473466807704
$ws.Cells.Item(2, 19).Value = 0.05156310915626622
$ws.Cells.Item(2, 20).Value = 0.05156310915626624
$ws.Cells.Item(3, 7).Value = 0.3484863333333333
$ws.Cells.Item(3, 8).Value = 1.045459
$ws.Cells.Item(3, 9).Value = 0.08515845388213966
$ws.Cells.Item(3, 10).Value = 0.08515845388213968
$ws.Cells.Item(3, 15).Value = 0.2540955070726236
$ws.Cells.Item(3, 16).Value = 0.2540955070726236
$ws.Cells.Item(3, 17).Value = 0.488352717162
$ws.Cells.Item(3, 18).Value = 4.395174454458
$ws.Cells.Item(3, 19).Value = 0.02163838052070291
$ws.Cells.Item(3, 20).Value = 0.02163838052070291
$ws.Cells.Item(4, 7).Value = 0.3484863333333333
$ws.Cells.Item(4, 8).Value = 1.045459
$ws.Cells.Item(4, 9).Value = 0.08515845388213966
$ws.Cells.Item(4, 10).Value = 0.08515845388213968
$ws.Cells.Item(4, 11).Value = 2.0
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1338136666666667
$ws.Cells.Item(4, 14).Value = 0.401441
$ws.Cells.Item(4, 15).Value = 0.02426328499787613
$ws.Cells.Item(4, 16).Value = 0.02426328499787612
$ws.Cells.Item(4, 17).Value = 0.04663223404655556
$ws.Cells.Item(4, 18).Value = 0.419690106419
$ws.Cells.Item(4, 19).Value = 0.002066223836520845
$ws.Cells.Item(4, 20).Value = 0.002066223836520845
$ws.Cells.Item(5, 7).Value = 0.3484863333333333
$ws.Cells.Item(5, 8).Value = 1.045459
$ws.Cells.Item(5, 9).Value = 0.08515845388213966
$ws.Cells.Item(5, 10).Value = 0.08515845388213968
$ws.Cells.Item(5, 13).Value = 0.6405483333333334
$ws.Cells.Item(5, 14).Value = 1.921645
$ws.Cells.Item(5, 15).Value = 0.11614513789011
$ws.Cells.Item(5, 16).Value = 0.11614513789011
$ws.Cells.Item(5, 17).Value = 0.2232223400061111
$ws.Cells.Item(5, 18).Value = 2.009001060055
$ws.Cells.Item(5, 19).Value = 0.009890740368649687
$ws.Cells.Item(5, 20).Value = 0.009890740368649687
$ws.Cells.Item(6, 9).Value = 0.6800250264078943
$ws.Cells.Item(6, 10).Value = 0.6800250264078944
$ws.Cells.Item(6, 13).Value = 3.339352
$ws.Cells.Item(6, 14).Value = 10.018056
$ws.Cells.Item(6, 15).Value = 0.6054960700393903
$ws.Cells.Item(6, 16).Value = 0.6054960700393903
$ws.Cells.Item(6, 17).Value = 9.292767668594667
$ws.Cells.Item(6, 18).Value = 83.634909017352
$ws.Cells.Item(6, 19).Value = 0.4117524810184126
$ws.Cells.Item(6, 20).Value = 0.4117524810184127
$ws.Cells.Item(7, 9).Value = 0.6800250264078943
$ws.Cells.Item(7, 10).Value = 0.6800250264078944
$ws.Cells.Item(7, 15).Value = 0.2540955070726236
$ws.Cells.Item(7, 16).Value = 0.2540955070726236
$ws.Cells.Item(7, 19).Value = 0.1727913039071882
$ws.Cells.Item(7, 20).Value = 0.1727913039071882
$ws.Cells.Item(8, 9).Value = 0.6800250264078943
$ws.Cells.Item(8, 10).Value = 0.6800250264078944
$ws.Cells.Item(8, 11).Value = 2.0
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1338136666666667
$ws.Cells.Item(8, 14).Value = 0.401441
$ws.Cells.Item(8, 15).Value = 0.02426328499787613
$ws.Cells.Item(8, 16).Value = 0.02426328499787612
$ws.Cells.Item(8, 17).Value = 0.3723774298774445
$ws.Cells.Item(8, 18).Value = 3.351396868897
$ws.Cells.Item(8, 19).Value = 0.01649964102142298
$ws.Cells.Item(8, 20).Value = 0.01649964102142298
$ws.Cells.Item(9, 9).Value = 0.6800250264078943
$ws.Cells.Item(9, 10).Value = 0.6800250264078944
$ws.Cells.Item(9, 13).Value = 0.6405483333333334
$ws.Cells.Item(9, 14).Value = 1.921645
$ws.Cells.Item(9, 15).Value = 0.11614513789011
$ws.Cells.Item(9, 16).Value = 0.11614513789011
$ws.Cells.Item(9, 17).Value = 1.782521531773889
$ws.Cells.Item(9, 18).Value = 16.042693785965
$ws.Cells.Item(9, 19).Value = 0.0789816004608706
$ws.Cells.Item(9, 20).Value = 0.0789816004608706
$ws.Cells.Item(10, 5).Value = 2.0
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.5631773333333333
$ws.Cells.Item(10, 8).Value = 1.689532
$ws.Cells.Item(10, 9).Value = 0.137621784215736
$ws.Cells.Item(10, 10).Value = 0.1376217842157361
$ws.Cells.Item(10, 13).Value = 3.339352
$ws.Cells.Item(10, 14).Value = 10.018056
$ws.Cells.Item(10, 15).Value = 0.6054960700393903
$ws.Cells.Item(10, 16).Value = 0.6054960700393903
$ws.Cells.Item(10, 17).Value = 1.880647354421333
$ws.Cells.Item(10, 18).Value = 16.925826189792
$ws.Cells.Item(10, 19).Value = 0.08332944949443717
$ws.Cells.Item(10, 20).Value = 0.08332944949443719
$ws.Cells.Item(11, 5).Value = 2.0
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.5631773333333333
$ws.Cells.Item(11, 8).Value = 1.689532
$ws.Cells.Item(11, 9).Value = 0.137621784215736
$ws.Cells.Item(11, 10).Value = 0.1376217842157361
$ws.Cells.Item(11, 15).Value = 0.2540955070726236
$ws.Cells.Item(11, 16).Value = 0.2540955070726236
$ws.Cells.Item(11, 17).Value = 0.789210808776
$ws.Cells.Item(11, 18).Value = 7.102897278984001
$ws.Cells.Item(11, 19).Value = 0.03496907704453664
$ws.Cells.Item(11, 20).Value = 0.03496907704453664
$ws.Cells.Item(12, 5).Value = 2.0
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.5631773333333333
$ws.Cells.Item(12, 8).Value = 1.689532
$ws.Cells.Item(12, 9).Value = 0.137621784215736
$ws.Cells.Item(12, 10).Value = 0.1376217842157361
$ws.Cells.Item(12, 11).Value = 2.0
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.1338136666666667
$ws.Cells.Item(12, 14).Value = 0.401441
$ws.Cells.Item(12, 15).Value = 0.02426328499787613
$ws.Cells.Item(12, 16).Value = 0.02426328499787612
$ws.Cells.Item(12, 17).Value = 0.0753608239568889
$ws.Cells.Item(12, 18).Value = 0.6782474156120001
$ws.Cells.Item(12, 19).Value = 0.003339156572342614
$ws.Cells.Item(12, 20).Value = 0.003339156572342614
$ws.Cells.Item(13, 5).Value = 2.0
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.5631773333333333
$ws.Cells.Item(13, 8).Value = 1.689532
$ws.Cells.Item(13, 9).Value = 0.137621784215736
$ws.Cells.Item(13, 10).Value = 0.1376217842157361
$ws.Cells.Item(13, 13).Value = 0.6405483333333334
$ws.Cells.Item(13, 14).Value = 1.921645
$ws.Cells.Item(13, 15).Value = 0.11614513789011
$ws.Cells.Item(13, 16).Value = 0.11614513789011
$ws.Cells.Item(13, 17).Value = 0.3607423022377778
$ws.Cells.Item(13, 18).Value = 3.24668072014
$ws.Cells.Item(13, 19).Value = 0.01598410110441963
$ws.Cells.Item(13, 20).Value = 0.01598410110441963
$ws.Cells.Item(14, 5).Value = 2.0
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.2482816666666667
$ws.Cells.Item(14, 8).Value = 0.744845
$ws.Cells.Item(14, 9).Value = 0.06067177056378329
$ws.Cells.Item(14, 10).Value = 0.0606717705637833
$ws.Cells.Item(14, 13).Value = 3.339352
$ws.Cells.Item(14, 14).Value = 10.018056
$ws.Cells.Item(14, 15).Value = 0.6054960700393903
$ws.Cells.Item(14, 16).Value = 0.6054960700393903
$ws.Cells.Item(14, 17).Value = 0.8290998801466667
$ws.Cells.Item(14, 18).Value = 7.46189892132
$ws.Cells.Item(14, 19).Value = 0.03673651863870235
$ws.Cells.Item(14, 20).Value = 0.03673651863870236
$ws.Cells.Item(15, 5).Value = 2.0
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.2482816666666667
$ws.Cells.Item(15, 8).Value = 0.744845
$ws.Cells.Item(15, 9).Value = 0.06067177056378329
$ws.Cells.Item(15, 10).Value = 0.0606717705637833
$ws.Cells.Item(15, 15).Value = 0.2540955070726236
$ws.Cells.Item(15, 16).Value = 0.2540955070726236
$ws.Cells.Item(15, 17).Value = 0.34793050671
$ws.Cells.Item(15, 18).Value = 3.13137456039
$ws.Cells.Item(15, 19).Value = 0.0154164243063984
$ws.Cells.Item(15, 20).Value = 0.0154164243063984
$ws.Cells.Item(16, 5).Value = 2.0
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.2482816666666667
$ws.Cells.Item(16, 8).Value = 0.744845
$ws.Cells.Item(16, 9).Value = 0.06067177056378329
$ws.Cells.Item(16, 10).Value = 0.0606717705637833
$ws.Cells.Item(16, 11).Value = 2.0
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.1338136666666667
$ws.Cells.Item(16, 14).Value = 0.401441
$ws.Cells.Item(16, 15).Value = 0.02426328499787613
$ws.Cells.Item(16, 16).Value = 0.02426328499787612
$ws.Cells.Item(16, 17).Value = 0.03322348018277778
$ws.Cells.Item(16, 18).Value = 0.299011321645
$ws.Cells.Item(16, 19).Value = 0.001472096460514825
$ws.Cells.Item(16, 20).Value = 0.001472096460514825
$ws.Cells.Item(17, 5).Value = 2.0
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.2482816666666667
$ws.Cells.Item(17, 8).Value = 0.744845
$ws.Cells.Item(17, 9).Value = 0.06067177056378329
$ws.Cells.Item(17, 10).Value = 0.0606717705637833
$ws.Cells.Item(17, 13).Value = 0.6405483333333334
$ws.Cells.Item(17, 14).Value = 1.921645
$ws.Cells.Item(17, 15).Value = 0.11614513789011
$ws.Cells.Item(17, 16).Value = 0.11614513789011
$ws.Cells.Item(17, 17).Value = 0.1590364077805556
$ws.Cells.Item(17, 18).Value = 1.431327670025
$ws.Cells.Item(17, 19).Value = 0.007046731158167729
$ws.Cells.Item(17, 20).Value = 0.00704673115816773
$ws.Cells.Item(18, 7).Value = 0.1494596666666667
$ws.Cells.Item(18, 8).Value = 0.448379
$ws.Cells.Item(18, 9).Value = 0.03652296493044672
$ws.Cells.Item(18, 10).Value = 0.03652296493044672
$ws.Cells.Item(18, 13).Value = 3.339352
$ws.Cells.Item(18, 14).Value = 10.018056
$ws.Cells.Item(18, 15).Value = 0.6054960700393903
$ws.Cells.Item(18, 16).Value = 0.6054960700393903
$ws.Cells.Item(18, 17).Value = 0.4990984368026667
$ws.Cells.Item(18, 18).Value = 4.491885931224
$ws.Cells.Item(18, 19).Value = 0.02211451173157197
$ws.Cells.Item(18, 20).Value = 0.02211451173157197
$ws.Cells.Item(19, 7).Value = 0.1494596666666667
$ws.Cells.Item(19, 8).Value = 0.448379
$ws.Cells.Item(19, 9).Value = 0.03652296493044672
$ws.Cells.Item(19, 10).Value = 0.03652296493044672
$ws.Cells.Item(19, 15).Value = 0.2540955070726236
$ws.Cells.Item(19, 16).Value = 0.2540955070726236
$ws.Cells.Item(19, 17).Value = 0.209445901722
$ws.Cells.Item(19, 18).Value = 1.885013115498
$ws.Cells.Item(19, 19).Value = 0.00928032129379751
$ws.Cells.Item(19, 20).Value = 0.009280321293797508
$ws.Cells.Item(20, 7).Value = 0.1494596666666667
$ws.Cells.Item(20, 8).Value = 0.448379
$ws.Cells.Item(20, 9).Value = 0.03652296493044672
$ws.Cells.Item(20, 10).Value = 0.03652296493044672
$ws.Cells.Item(20, 11).Value = 2.0
$ws.Cells.Item(20, 12).Value = 0.6666666666666666
$ws.Cells.Item(20, 13).Value = 0.1338136666666667
$ws.Cells.Item(20, 14).Value = 0.401441
$ws.Cells.Item(20, 15).Value = 0.02426328499787613
$ws.Cells.Item(20, 16).Value = 0.02426328499787612
$ws.Cells.Item(20, 17).Value = 0.01999974601544445
$ws.Cells.Item(20, 18).Value = 0.179997714139
$ws.Cells.Item(20, 19).Value = 0.0008861671070748638
$ws.Cells.Item(20, 20).Value = 0.0008861671070748637
$ws.Cells.Item(21, 7).Value = 0.1494596666666667
$ws.Cells.Item(21, 8).Value = 0.448379
$ws.Cells.Item(21, 9).Value = 0.03652296493044672
$ws.Cells.Item(21, 10).Value = 0.03652296493044672
$ws.Cells.Item(21, 13).Value = 0.6405483333333334
$ws.Cells.Item(21, 14).Value = 1.921645
$ws.Cells.Item(21, 15).Value = 0.11614513789011
$ws.Cells.Item(21, 16).Value = 0.11614513789011
$ws.Cells.Item(21, 17).Value = 0.09573614038388889
$ws.Cells.Item(21, 18).Value = 0.861625263455
$ws.Cells.Item(21, 19).Value = 0.004241964798002388
$ws.Cells.Item(21, 20).Value = 0.004241964798002387
